$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 768, shifting existing rows 768-818 down to 772-822
$ws.Range("A768:T771").EntireRow.Insert()

# Populate the 4 newly inserted rows with the new weekly price records
$newRows = @(
    @{ Row=768; D=44746; K="Fukumoto";   L="Primera"; M=143; N=4000; O=4500; P=4238; S=326 },
    @{ Row=769; D=44746; K="Fukumoto";   L="Segunda"; M=165; N=3000; O=3500; P=3258; S=251 },
    @{ Row=770; D=44746; K="Navel Late"; L="Primera"; M=187; N=4000; O=4500; P=4262; S=328 },
    @{ Row=771; D=44746; K="Navel Late"; L="Segunda"; M=153; N=3000; O=3500; P=3255; S=250 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102005
    $ws.Cells.Item($row, 10).Value = "Naranja"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/malla 13 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Quillota"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 13
}

Write-Host "Inserted and populated rows"
